# Auto-derived edits applying the Seraph_Profits.xlsx diff
# Updates price/profit columns (H:N) for specific Leve rows across several sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 38
$ws.Range("H38").Value = 779
$ws.Range("I38").Value = 271.16666
$ws.Range("J38").Value = 1214.2858
$ws.Range("K38").Value = 813.4999799999999
$ws.Range("L38").Value = 3642.8574
$ws.Range("M38").Value = -441.4999799999999
$ws.Range("N38").Value = -4386.857400000001
# Row 43
$ws.Range("H43").Value = 12083.083
$ws.Range("I43").Value = 6875
$ws.Range("J43").Value = 14687.125
$ws.Range("K43").Value = 6875
$ws.Range("L43").Value = 14687.125
$ws.Range("M43").Value = -6806
$ws.Range("N43").Value = -14825.125
# Row 134
$ws.Range("H134").Value = 64000
$ws.Range("J134").Value = 64000
$ws.Range("L134").Value = 64000
$ws.Range("N134").Value = -74140
# Row 141
$ws.Range("H141").Value = 10495.667
$ws.Range("I141").Value = 5993.3335
$ws.Range("K141").Value = 17980.0005
$ws.Range("M141").Value = -12800.0005

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 6069.7144
$ws.Range("I61").Value = 6069.7144
$ws.Range("K61").Value = 6069.7144
$ws.Range("M61").Value = -5857.7144
# Row 130
$ws.Range("H130").Value = 44714
$ws.Range("J130").Value = 44714
$ws.Range("L130").Value = 44714
$ws.Range("N130").Value = -54754
# Row 136
$ws.Range("H136").Value = 6069.7144
$ws.Range("I136").Value = 6069.7144
$ws.Range("K136").Value = 18209.1432
$ws.Range("M136").Value = -15659.1432

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 3043.6128
$ws.Range("I105").Value = 2726.5
$ws.Range("K105").Value = 2726.5
$ws.Range("M105").Value = -979.5
# Row 134
$ws.Range("H134").Value = 424.5
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null

$ws = $wb.Worksheets.Item("CRP")
# Row 3
$ws.Range("H3").Value = 20471.875
$ws.Range("I3").Value = 18946
$ws.Range("J3").Value = 21997.75
$ws.Range("K3").Value = 18946
$ws.Range("L3").Value = 21997.75
$ws.Range("M3").Value = -18833
$ws.Range("N3").Value = -22223.75
# Row 13
$ws.Range("H13").Value = 1512
$ws.Range("I13").Value = 765
$ws.Range("J13").Value = 1885.5
$ws.Range("K13").Value = 765
$ws.Range("L13").Value = 1885.5
$ws.Range("M13").Value = -626
$ws.Range("N13").Value = -2163.5
# Row 15
$ws.Range("H15").Value = 4578.75
$ws.Range("I15").Value = 7407.5
$ws.Range("K15").Value = 7407.5
$ws.Range("M15").Value = -7237.5
# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = $null
$ws.Range("N45").Value = 0
# Row 105
$ws.Range("H105").Value = 1966.4166
$ws.Range("I105").Value = 971
$ws.Range("K105").Value = 971
$ws.Range("M105").Value = 776
# Row 107
$ws.Range("H107").Value = 1340.9474
$ws.Range("I107").Value = 1094.4615
$ws.Range("J107").Value = 1875
$ws.Range("K107").Value = 1094.4615
$ws.Range("L107").Value = 1875
$ws.Range("M107").Value = 825.5385000000001
$ws.Range("N107").Value = -5715

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 14285884
$ws.Range("I7").Value = 25000174
$ws.Range("J7").Value = 163.66667
$ws.Range("K7").Value = 75000522
$ws.Range("L7").Value = 491.00001
$ws.Range("M7").Value = -75000410
$ws.Range("N7").Value = -715.00001
# Row 140
$ws.Range("H140").Value = 3332.889
$ws.Range("I140").Value = 3124.5
$ws.Range("K140").Value = 9373.5
$ws.Range("M140").Value = -4193.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3756.5
$ws.Range("I80").Value = 2916.6667
$ws.Range("J80").Value = 3985.5454
$ws.Range("K80").Value = 2916.6667
$ws.Range("L80").Value = 3985.5454
$ws.Range("M80").Value = -1918.6667
$ws.Range("N80").Value = -5981.5454
# Row 83
$ws.Range("H83").Value = 3756.5
$ws.Range("I83").Value = 2916.6667
$ws.Range("J83").Value = 3985.5454
$ws.Range("K83").Value = 14583.3335
$ws.Range("L83").Value = 19927.727
$ws.Range("M83").Value = -9591.333500000001
$ws.Range("N83").Value = -29911.727
# Row 132
$ws.Range("H132").Value = 1458.8125
$ws.Range("I132").Value = 1238.7142
$ws.Range("K132").Value = 3716.1426
$ws.Range("M132").Value = -1186.1426

$ws = $wb.Worksheets.Item("LTW")
# Row 2
$ws.Range("H2").Value = 10010001
$ws.Range("I2").Value = 20000000
$ws.Range("J2").Value = 20002
$ws.Range("K2").Value = 20000000
$ws.Range("L2").Value = 20002
$ws.Range("M2").Value = -19999888
$ws.Range("N2").Value = -20226
# Row 34
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = $null
$ws.Range("M34").Value = $null
$ws.Range("N34").Value = 0
# Row 55
$ws.Range("H55").Value = 1009.2857
$ws.Range("I55").Value = 824.25
$ws.Range("K55").Value = 824.25
$ws.Range("M55").Value = -651.25
# Row 68
$ws.Range("H68").Value = 4366.3335
$ws.Range("I68").Value = 4286.75
$ws.Range("J68").Value = 5003
$ws.Range("K68").Value = 4286.75
$ws.Range("L68").Value = 5003
$ws.Range("M68").Value = -3537.75
$ws.Range("N68").Value = -6501
# Row 71
$ws.Range("H71").Value = 4366.3335
$ws.Range("I71").Value = 4286.75
$ws.Range("J71").Value = 5003
$ws.Range("K71").Value = 21433.75
$ws.Range("L71").Value = 25015
$ws.Range("M71").Value = -17689.75
$ws.Range("N71").Value = -32503
# Row 136
$ws.Range("H136").Value = 6530.4
$ws.Range("I136").Value = 6362
$ws.Range("K136").Value = 19086
$ws.Range("M136").Value = -16536

$ws = $wb.Worksheets.Item("WVR")
# Row 12
$ws.Range("H12").Value = 4279.6
$ws.Range("I12").Value = 4249.5
$ws.Range("J12").Value = 4299.6665
$ws.Range("K12").Value = 4249.5
$ws.Range("L12").Value = 4299.6665
$ws.Range("M12").Value = -4107.5
$ws.Range("N12").Value = -4583.6665
# Row 33
$ws.Range("H33").Value = 40333
$ws.Range("J33").Value = 45000
$ws.Range("L33").Value = 45000
$ws.Range("N33").Value = -45500
# Row 36
$ws.Range("H36").Value = 40333
$ws.Range("J36").Value = 45000
$ws.Range("L36").Value = 45000
$ws.Range("N36").Value = -45500
# Row 40
$ws.Range("H40").Value = 38083.168
$ws.Range("J40").Value = 38299.8
$ws.Range("L40").Value = 38299.8
$ws.Range("N40").Value = -38597.8
# Row 56
$ws.Range("H56").Value = 54975
$ws.Range("I56").Value = 54975
$ws.Range("K56").Value = 54975
$ws.Range("M56").Value = -54261
# Row 62
$ws.Range("H62").Value = 4666.0557
$ws.Range("J62").Value = 4799.3335
$ws.Range("L62").Value = 4799.3335
$ws.Range("N62").Value = -6047.3335
# Row 65
$ws.Range("H65").Value = 4666.0557
$ws.Range("J65").Value = 4799.3335
$ws.Range("L65").Value = 23996.6675
$ws.Range("N65").Value = -30236.6675
# Row 76
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = $null
$ws.Range("N76").Value = 0
# Row 79
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = $null
$ws.Range("N79").Value = 0
# Row 81
$ws.Range("H81").Value = 4960.643
$ws.Range("I81").Value = 5791.6665
$ws.Range("J81").Value = 4337.375
$ws.Range("K81").Value = 11583.333
$ws.Range("L81").Value = 8674.75
$ws.Range("M81").Value = -10522.333
$ws.Range("N81").Value = -10796.75
# Row 84
$ws.Range("H84").Value = 4960.643
$ws.Range("I84").Value = 5791.6665
$ws.Range("J84").Value = 4337.375
$ws.Range("K84").Value = 57916.665
$ws.Range("L84").Value = 43373.75
$ws.Range("M84").Value = -52612.665
$ws.Range("N84").Value = -53981.75
# Row 93
$ws.Range("H93").Value = 50000
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("M93").Value = -54992
# Row 113
$ws.Range("H113").Value = 290.27274
$ws.Range("I113").Value = 262
$ws.Range("J113").Value = 365.66666
$ws.Range("K113").Value = 786
$ws.Range("L113").Value = 1096.99998
$ws.Range("M113").Value = 1384
$ws.Range("N113").Value = -5436.999980000001
# Row 122
$ws.Range("H122").Value = 2445.7
$ws.Range("I122").Value = 2682.4285
$ws.Range("J122").Value = 1893.3334
$ws.Range("K122").Value = 8047.2855
$ws.Range("L122").Value = 5680.0002
$ws.Range("M122").Value = -5597.2855
$ws.Range("N122").Value = -10580.0002
# Row 136
$ws.Range("H136").Value = 892.43475
$ws.Range("I136").Value = 667.9048
$ws.Range("K136").Value = 2003.7144
$ws.Range("M136").Value = 546.2855999999999

Write-Output "Applied Seraph_Profits updates."
